$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# --- Row 17: the measurement columns (B:K) and the Event column (M) were left
# blank when this row was first created; this edit fills them in with the
# literal placeholder text "nan" (the rest of the sheet already uses "nan" as
# the textual stand-in for "no value").
$ws.Range("B17:K17").Value = "nan"
$ws.Range("M17").Value = "nan"

# --- Row 18: brand new service-history entry for Card22 ("إضافة حدث جديد في Card22").
# Column A keeps the "card" number as text (matches every other row in the
# column), the measurement columns B:K have no data for this event, and the
# new service record itself lives in Date / Event / Correction / Serviced by.
$ws.Range("A18").Value = "'22"
$ws.Range("L18").Value = "14\5\2025"
$ws.Range("M18").Value = "629.7 t"
$ws.Range("N18").Value = "تم عمل صيانه وتغيير الجرائد الاماميه ومعايره المكنه (1_2_5_7_8)"
$ws.Range("O18").Value = "الخبير"
